$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update summary figures at the top of the statement ---
# Valor Mora (total overdue amount)
$ws.Range("E11").Value = 113880
# Cant. Trabajadores (worker count) - now 2 workers
$ws.Range("C13").Value = 2
# Cant. Periodos (period count) - now only 1 period
$ws.Range("F13").Value = 1

# --- Rebuild the detail table ---
# Previously the table listed the SAME worker (Santiago) across 5 duplicate
# period rows (16-20, periods 2507..2503). Row 20 carried the special
# "closing" border style for the bottom of the table.
#
# New layout: row 16 becomes the new worker (Kenia Liz Diaz Perez) for
# period 2508, and row 17 keeps the existing worker (Santiago) but updated
# to period 2508, now taking on the table's closing border style (which
# used to belong to row 20). The old rows 18-20 are removed.

# Move row 20's formatting/values (the closing-border row) up onto row 17.
$ws.Range("B20:J20").Copy($ws.Range("B17:J17"))

# Remove the now-duplicate rows 18, 19, 20.
$ws.Rows("18:20").Delete()

# Update worker data:
# Row 16: replace with the new worker's info for period 2508.
$ws.Range("C16").Value = "1143413842"
$ws.Range("D16").Value = "KENIA LIZ DIAZ PEREZ"
$ws.Range("E16").Value = "2508"

# Row 17: existing worker (Santiago), period updated to 2508.
$ws.Range("E17").Value = "2508"

Write-Host "Done."
